$wb = $excel.ActiveWorkbook

# Work on the "line_imp" sheet (sheet2) - add two new header columns.
$ws = $wb.Worksheets.Item("line_imp")

$ws.Range("F1").Value = "t_x"
$ws.Range("G1").Value = "t_a"

# Select G1 and make this sheet the active one (matches the diff's selection/activeTab changes)
$ws.Activate()
$ws.Range("G1").Select()
